$wb = $excel.ActiveWorkbook

# --- Step 1: rename all sheets to unique temporary names to avoid name collisions ---
$wb.Worksheets.Item(1).Name = "__tmp1__"
$wb.Worksheets.Item(2).Name = "__tmp2__"
$wb.Worksheets.Item(3).Name = "__tmp3__"
$wb.Worksheets.Item(4).Name = "__tmp4__"
$wb.Worksheets.Item(5).Name = "__tmp5__"
$wb.Worksheets.Item(6).Name = "__tmp6__"
$wb.Worksheets.Item(7).Name = "__tmp7__"
$wb.Worksheets.Item(8).Name = "__tmp8__"
$wb.Worksheets.Item(9).Name = "__tmp9__"

# --- Step 2: rename to final target names (tab order / r:id order unchanged) ---
$wb.Worksheets.Item(1).Name = "summ7"
$wb.Worksheets.Item(2).Name = "summ10"
$wb.Worksheets.Item(3).Name = "summ2"
$wb.Worksheets.Item(4).Name = "summ1"
$wb.Worksheets.Item(5).Name = "summ9"
$wb.Worksheets.Item(6).Name = "summ3"
$wb.Worksheets.Item(7).Name = "summ5"
$wb.Worksheets.Item(8).Name = "summ4"
$wb.Worksheets.Item(9).Name = "summ0"

# --- Update data for sheet (physical position 1): summ4 -> summ7 ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = [double]"-0.4243636629232816"
$ws.Range("C2").Value = [double]"0.798219854410027"
$ws.Range("B3").Value = [double]"0.8742129440741078"
$ws.Range("C3").Value = [double]"0.585172252358718"
$ws.Range("B4").Value = [double]"1.060859301216261"
$ws.Range("C4").Value = [double]"0.5096033400296811"
$ws.Range("B5").Value = [double]"-0.3479172045388123"
$ws.Range("C5").Value = [double]"0.8282537980304184"
$ws.Range("B6").Value = [double]"-0.1217643990686017"
$ws.Range("C6").Value = [double]"0.9394916577723701"
$ws.Range("B7").Value = [double]"0.2002211151863326"
$ws.Range("C7").Value = [double]"0.9009280793337782"
$ws.Range("B8").Value = [double]"-0.08959795388004786"
$ws.Range("C8").Value = [double]"0.2223632416295082"
$ws.Range("B9").Value = [double]"0.0005104152727709657"
$ws.Range("C9").Value = [double]"1.511260390179557e-63"
$ws.Range("B10").Value = [double]"0.009051788646237645"
$ws.Range("C10").Value = [double]"0.0008062653047823078"
$ws.Range("B11").Value = [double]"-0.1820040854889209"
$ws.Range("C11").Value = [double]"0.01044012607352036"
$ws.Range("B12").Value = [double]"0.81021682738922"
$ws.Range("C12").Value = [double]"3.919416002181223e-14"
$ws.Range("B13").Value = [double]"0.2591170652767112"
$ws.Range("C13").Value = [double]"0.06037004442028996"
$ws.Range("B14").Value = [double]"-5.134292847833611e-05"
$ws.Range("C14").Value = [double]"0.0009518291279877662"
$ws.Range("B15").Value = [double]"-2.064388343054366e-08"
$ws.Range("C15").Value = [double]"0.1869271944639156"
$ws.Range("B16").Value = [double]"0.01127095230283927"
$ws.Range("C16").Value = [double]"0.6541483380468234"
$ws.Range("B17").Value = [double]"0.0869085651529721"
$ws.Range("C17").Value = [double]"3.039763814956007e-07"
$ws.Range("B18").Value = [double]"-1.406667345217116"
$ws.Range("C18").Value = [double]"7.877374563171531e-06"
$ws.Range("B19").Value = [double]"-0.009842140736899"
$ws.Range("C19").Value = [double]"0.01384710392559362"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.003826690173774339"
$ws.Range("C20").Value = [double]"0.1037077996943826"
$ws.Range("B21").Value = [double]"0.4985993148202589"
$ws.Range("C21").Value = [double]"0.1761829466310516"
$ws.Range("B22").Value = [double]"0.002199501766126337"
$ws.Range("C22").Value = [double]"0.9945415429126517"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 2): summ1 -> summ10 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = [double]"-1.086151981348478"
$ws.Range("C2").Value = [double]"0.4488652153402993"
$ws.Range("B3").Value = [double]"1.323823823383789"
$ws.Range("C3").Value = [double]"0.3332555044268218"
$ws.Range("B4").Value = [double]"1.474533224988986"
$ws.Range("C4").Value = [double]"0.283864005387796"
$ws.Range("B5").Value = [double]"0.1511349342597691"
$ws.Range("C5").Value = [double]"0.9121807360040631"
$ws.Range("B6").Value = [double]"0.4291221183986465"
$ws.Range("C6").Value = [double]"0.7542421304402194"
$ws.Range("B7").Value = [double]"0.6393254815791174"
$ws.Range("C7").Value = [double]"0.6421837821349758"
$ws.Range("B8").Value = [double]"-0.04241568976596242"
$ws.Range("C8").Value = [double]"0.5486999321671051"
$ws.Range("B9").Value = [double]"0.0004978006249782331"
$ws.Range("C9").Value = [double]"2.787242841734438e-61"
$ws.Range("B10").Value = [double]"0.007729574933440251"
$ws.Range("C10").Value = [double]"0.00390099142593696"
$ws.Range("B11").Value = [double]"-0.09988879981345412"
$ws.Range("C11").Value = [double]"0.1564458243487198"
$ws.Range("B12").Value = [double]"0.8477981781031829"
$ws.Range("C12").Value = [double]"1.544505875745745e-15"
$ws.Range("B13").Value = [double]"0.3559703190049269"
$ws.Range("C13").Value = [double]"0.009161954817653917"
$ws.Range("B14").Value = [double]"-4.725100697854248e-05"
$ws.Range("C14").Value = [double]"0.002300759760269056"
$ws.Range("B15").Value = [double]"-2.27117046914282e-08"
$ws.Range("C15").Value = [double]"0.1514210746984785"
$ws.Range("B16").Value = [double]"0.002777075255159134"
$ws.Range("C16").Value = [double]"0.911681932244361"
$ws.Range("B17").Value = [double]"0.09546780677604239"
$ws.Range("C17").Value = [double]"1.45571291798196e-08"
$ws.Range("B18").Value = [double]"-1.72184546904808"
$ws.Range("C18").Value = [double]"3.249982572740163e-08"
$ws.Range("B19").Value = [double]"-0.007675436958777456"
$ws.Range("C19").Value = [double]"0.05370378281814721"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.003534518239995371"
$ws.Range("C20").Value = [double]"0.1321117736652766"
$ws.Range("B21").Value = [double]"0.3570116917130017"
$ws.Range("C21").Value = [double]"0.3305289497421604"
$ws.Range("B22").Value = [double]"0.0648786032437135"
$ws.Range("C22").Value = [double]"0.8414734615967516"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 3): summ10 -> summ2 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = [double]"-0.8660553420481917"
$ws.Range("C2").Value = [double]"0.5508668110950015"
$ws.Range("B3").Value = [double]"1.366926314605128"
$ws.Range("C3").Value = [double]"0.3237086872539922"
$ws.Range("B4").Value = [double]"1.535913958071382"
$ws.Range("C4").Value = [double]"0.2701442208890351"
$ws.Range("B5").Value = [double]"0.114951934825112"
$ws.Range("C5").Value = [double]"0.9339628789923596"
$ws.Range("B6").Value = [double]"0.373173130632831"
$ws.Range("C6").Value = [double]"0.7879945260634684"
$ws.Range("B7").Value = [double]"0.8120746977602343"
$ws.Range("C7").Value = [double]"0.5599822715073262"
$ws.Range("B8").Value = [double]"-0.06553642047331024"
$ws.Range("C8").Value = [double]"0.359842129898239"
$ws.Range("B9").Value = [double]"0.0004930186106153967"
$ws.Range("C9").Value = [double]"4.148340320063387e-60"
$ws.Range("B10").Value = [double]"0.00976492934945531"
$ws.Range("C10").Value = [double]"0.0003034133761173874"
$ws.Range("B11").Value = [double]"-0.1286480411705949"
$ws.Range("C11").Value = [double]"0.07090596991298713"
$ws.Range("B12").Value = [double]"0.837653359128266"
$ws.Range("C12").Value = [double]"8.661206785172269e-15"
$ws.Range("B13").Value = [double]"0.3149954625500732"
$ws.Range("C13").Value = [double]"0.0236198131239718"
$ws.Range("B14").Value = [double]"-4.805483198402628e-05"
$ws.Range("C14").Value = [double]"0.002032748421626249"
$ws.Range("B15").Value = [double]"-2.531438846596676e-08"
$ws.Range("C15").Value = [double]"0.1065099928460105"
$ws.Range("B16").Value = [double]"0.01718563744415141"
$ws.Range("C16").Value = [double]"0.4950731746959248"
$ws.Range("B17").Value = [double]"0.08240349993483097"
$ws.Range("C17").Value = [double]"1.177292598538041e-06"
$ws.Range("B18").Value = [double]"-1.455214051234775"
$ws.Range("C18").Value = [double]"3.376826608278778e-06"
$ws.Range("B19").Value = [double]"-0.01094892414356155"
$ws.Range("C19").Value = [double]"0.00619987426917606"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.004545495588290685"
$ws.Range("C20").Value = [double]"0.05570616320676525"
$ws.Range("B21").Value = [double]"0.4385968791566145"
$ws.Range("C21").Value = [double]"0.2343156626418368"
$ws.Range("B22").Value = [double]"0.007516626929020852"
$ws.Range("C22").Value = [double]"0.9815739921959711"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 4): summ8 -> summ1 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = [double]"-16.37959239501854"
$ws.Range("C2").Value = [double]"0.9961409113967011"
$ws.Range("B3").Value = [double]"16.84127318011505"
$ws.Range("C3").Value = [double]"0.9960321386604736"
$ws.Range("B4").Value = [double]"17.04223804963526"
$ws.Range("C4").Value = [double]"0.9959847910535579"
$ws.Range("B5").Value = [double]"15.58845700843931"
$ws.Range("C5").Value = [double]"0.9963273042512173"
$ws.Range("B6").Value = [double]"15.90343114818129"
$ws.Range("C6").Value = [double]"0.9962530957671077"
$ws.Range("B7").Value = [double]"16.23368137140249"
$ws.Range("C7").Value = [double]"0.996175288246744"
$ws.Range("B8").Value = [double]"-0.08503583589562548"
$ws.Range("C8").Value = [double]"0.2320272208688644"
$ws.Range("B9").Value = [double]"0.0005042967082253355"
$ws.Range("C9").Value = [double]"4.269796242762537e-63"
$ws.Range("B10").Value = [double]"0.008001241341365754"
$ws.Range("C10").Value = [double]"0.002949403643478083"
$ws.Range("B11").Value = [double]"-0.171609484969568"
$ws.Range("C11").Value = [double]"0.01584393652152533"
$ws.Range("B12").Value = [double]"0.7749128326173912"
$ws.Range("C12").Value = [double]"3.860051062166705e-13"
$ws.Range("B13").Value = [double]"0.2964306555508469"
$ws.Range("C13").Value = [double]"0.03212518493102407"
$ws.Range("B14").Value = [double]"-5.276943213411192e-05"
$ws.Range("C14").Value = [double]"0.0006914429599640659"
$ws.Range("B15").Value = [double]"-2.09505078127529e-08"
$ws.Range("C15").Value = [double]"0.1833923291946958"
$ws.Range("B16").Value = [double]"0.006512352765672753"
$ws.Range("C16").Value = [double]"0.7946170759410639"
$ws.Range("B17").Value = [double]"0.08868083110675073"
$ws.Range("C17").Value = [double]"1.224212784364523e-07"
$ws.Range("B18").Value = [double]"-1.377660778348481"
$ws.Range("C18").Value = [double]"1.117407457273061e-05"
$ws.Range("B19").Value = [double]"-0.009409075538244074"
$ws.Range("C19").Value = [double]"0.01884059727943451"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.003528771623753404"
$ws.Range("C20").Value = [double]"0.1342381366020181"
$ws.Range("B21").Value = [double]"0.4400350984497851"
$ws.Range("C21").Value = [double]"0.2332308879123858"
$ws.Range("B22").Value = [double]"0.102968343658936"
$ws.Range("C22").Value = [double]"0.7485048822480478"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 5): summ2 -> summ9 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = [double]"-0.503002900906211"
$ws.Range("C2").Value = [double]"0.760284824626409"
$ws.Range("B3").Value = [double]"0.8961744656371015"
$ws.Range("C3").Value = [double]"0.5733647203885981"
$ws.Range("B4").Value = [double]"1.034104556401019"
$ws.Range("C4").Value = [double]"0.5175710346815465"
$ws.Range("B5").Value = [double]"-0.4034586538857923"
$ws.Range("C5").Value = [double]"0.800094779174598"
$ws.Range("B6").Value = [double]"-0.06567907039793047"
$ws.Range("C6").Value = [double]"0.9671268327708611"
$ws.Range("B7").Value = [double]"0.3056381391848095"
$ws.Range("C7").Value = [double]"0.8483488187301718"
$ws.Range("B8").Value = [double]"-0.07845932081577503"
$ws.Range("C8").Value = [double]"0.2636150748040528"
$ws.Range("B9").Value = [double]"0.0004802184110018472"
$ws.Range("C9").Value = [double]"3.43462418292957e-58"
$ws.Range("B10").Value = [double]"0.00854645002163759"
$ws.Range("C10").Value = [double]"0.001586982072948138"
$ws.Range("B11").Value = [double]"-0.1294407635426284"
$ws.Range("C11").Value = [double]"0.0693478077626685"
$ws.Range("B12").Value = [double]"0.8978598032866361"
$ws.Range("C12").Value = [double]"3.167693158240906e-17"
$ws.Range("B13").Value = [double]"0.3737133780689666"
$ws.Range("C13").Value = [double]"0.006554467361421182"
$ws.Range("B14").Value = [double]"-4.688593232917798e-05"
$ws.Range("C14").Value = [double]"0.002315050491289776"
$ws.Range("B15").Value = [double]"-2.055764772215072e-08"
$ws.Range("C15").Value = [double]"0.1907120645841294"
$ws.Range("B16").Value = [double]"0.007072265604767533"
$ws.Range("C16").Value = [double]"0.7785513113387559"
$ws.Range("B17").Value = [double]"0.09298785337521767"
$ws.Range("C17").Value = [double]"4.745190388367354e-08"
$ws.Range("B18").Value = [double]"-1.645768531729923"
$ws.Range("C18").Value = [double]"1.821970686562076e-07"
$ws.Range("B19").Value = [double]"-0.007446639017054751"
$ws.Range("C19").Value = [double]"0.06129835174022009"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.003725464171268802"
$ws.Range("C20").Value = [double]"0.1133034113464239"
$ws.Range("B21").Value = [double]"0.1679636901280385"
$ws.Range("C21").Value = [double]"0.6463549330166222"
$ws.Range("B22").Value = [double]"0.09475800325055575"
$ws.Range("C22").Value = [double]"0.770964359651146"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 6): summ6 -> summ3 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = [double]"-1.042124414910794"
$ws.Range("C2").Value = [double]"0.4667862098903872"
$ws.Range("B3").Value = [double]"1.390750756456341"
$ws.Range("C3").Value = [double]"0.3078063647706812"
$ws.Range("B4").Value = [double]"1.533029310418023"
$ws.Range("C4").Value = [double]"0.2636925673396001"
$ws.Range("B5").Value = [double]"0.1433115735424497"
$ws.Range("C5").Value = [double]"0.9164390098365077"
$ws.Range("B6").Value = [double]"0.4675583636172128"
$ws.Range("C6").Value = [double]"0.732198347781343"
$ws.Range("B7").Value = [double]"0.6586682135401819"
$ws.Range("C7").Value = [double]"0.6310052641199202"
$ws.Range("B8").Value = [double]"-0.03162584963764769"
$ws.Range("C8").Value = [double]"0.6626549296580577"
$ws.Range("B9").Value = [double]"0.0004883630596948524"
$ws.Range("C9").Value = [double]"5.069628136861065e-58"
$ws.Range("B10").Value = [double]"0.0086287999942775"
$ws.Range("C10").Value = [double]"0.001510472414876151"
$ws.Range("B11").Value = [double]"-0.1374436517064476"
$ws.Range("C11").Value = [double]"0.05508663013878375"
$ws.Range("B12").Value = [double]"0.8087533977429244"
$ws.Range("C12").Value = [double]"6.304743583465699e-14"
$ws.Range("B13").Value = [double]"0.3027257128056383"
$ws.Range("C13").Value = [double]"0.02999822013374287"
$ws.Range("B14").Value = [double]"-4.895285292855531e-05"
$ws.Range("C14").Value = [double]"0.001886831612988988"
$ws.Range("B15").Value = [double]"-1.504467044116031e-08"
$ws.Range("C15").Value = [double]"0.3422604502537403"
$ws.Range("B16").Value = [double]"-0.003123853151686047"
$ws.Range("C16").Value = [double]"0.9020273151245845"
$ws.Range("B17").Value = [double]"0.08877857331800916"
$ws.Range("C17").Value = [double]"2.137598666067717e-07"
$ws.Range("B18").Value = [double]"-1.700914905725428"
$ws.Range("C18").Value = [double]"5.20022705580963e-08"
$ws.Range("B19").Value = [double]"-0.00900541835679566"
$ws.Range("C19").Value = [double]"0.02561846514044125"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.00334264756427975"
$ws.Range("C20").Value = [double]"0.1598932567762282"
$ws.Range("B21").Value = [double]"0.4625582730908633"
$ws.Range("C21").Value = [double]"0.2129664402511514"
$ws.Range("B22").Value = [double]"-0.0254435772773517"
$ws.Range("C22").Value = [double]"0.9375228375048478"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 7): summ3 -> summ5 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = [double]"-1.408776199468974"
$ws.Range("C2").Value = [double]"0.3278060677379981"
$ws.Range("B3").Value = [double]"1.323313492639765"
$ws.Range("C3").Value = [double]"0.3347276474170552"
$ws.Range("B4").Value = [double]"1.358223893234059"
$ws.Range("C4").Value = [double]"0.3249012243519389"
$ws.Range("B5").Value = [double]"0.1116438723344182"
$ws.Range("C5").Value = [double]"0.9352453765241923"
$ws.Range("B6").Value = [double]"0.420577471589437"
$ws.Range("C6").Value = [double]"0.7596113403083975"
$ws.Range("B7").Value = [double]"0.6593834985062097"
$ws.Range("C7").Value = [double]"0.6326667310909944"
$ws.Range("B8").Value = [double]"-0.01510939826574014"
$ws.Range("C8").Value = [double]"0.8351919483892996"
$ws.Range("B9").Value = [double]"0.0005082533540190904"
$ws.Range("C9").Value = [double]"5.523186940214031e-63"
$ws.Range("B10").Value = [double]"0.009318045851459596"
$ws.Range("C10").Value = [double]"0.0005769005418445238"
$ws.Range("B11").Value = [double]"-0.1471476242561881"
$ws.Range("C11").Value = [double]"0.03854957047195873"
$ws.Range("B12").Value = [double]"0.8906240382427038"
$ws.Range("C12").Value = [double]"9.378259424681416e-17"
$ws.Range("B13").Value = [double]"0.3233438893440145"
$ws.Range("C13").Value = [double]"0.01918283968089252"
$ws.Range("B14").Value = [double]"-5.277077938151113e-05"
$ws.Range("C14").Value = [double]"0.0007315648589865599"
$ws.Range("B15").Value = [double]"-2.601715004066161e-08"
$ws.Range("C15").Value = [double]"0.09847022484671535"
$ws.Range("B16").Value = [double]"-0.02225737430802336"
$ws.Range("C16").Value = [double]"0.3757770162721117"
$ws.Range("B17").Value = [double]"0.09714855839764384"
$ws.Range("C17").Value = [double]"9.802803748726448e-09"
$ws.Range("B18").Value = [double]"-1.521627617961461"
$ws.Range("C18").Value = [double]"1.303230018931014e-06"
$ws.Range("B19").Value = [double]"-0.005041836868323587"
$ws.Range("C19").Value = [double]"0.2082212985863918"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.002110454246813226"
$ws.Range("C20").Value = [double]"0.3741804811386299"
$ws.Range("B21").Value = [double]"0.3655542296460046"
$ws.Range("C21").Value = [double]"0.3231620229801482"
$ws.Range("B22").Value = [double]"0.1409544897241429"
$ws.Range("C22").Value = [double]"0.6605821573691142"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 8): summ0 -> summ4 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = [double]"-0.9153581566427168"
$ws.Range("C2").Value = [double]"0.5202284599491819"
$ws.Range("B3").Value = [double]"1.406386862406979"
$ws.Range("C3").Value = [double]"0.2991970288468754"
$ws.Range("B4").Value = [double]"1.542004685960225"
$ws.Range("C4").Value = [double]"0.2578551239398247"
$ws.Range("B5").Value = [double]"0.1514086499771623"
$ws.Range("C5").Value = [double]"0.9111581567845909"
$ws.Range("B6").Value = [double]"0.4736916651913417"
$ws.Range("C6").Value = [double]"0.7271063583215063"
$ws.Range("B7").Value = [double]"0.7038618960575357"
$ws.Range("C7").Value = [double]"0.6055185249715472"
$ws.Range("B8").Value = [double]"-0.03607724471865951"
$ws.Range("C8").Value = [double]"0.6177732058886267"
$ws.Range("B9").Value = [double]"0.0005027049049543321"
$ws.Range("C9").Value = [double]"5.877139484504336e-61"
$ws.Range("B10").Value = [double]"0.006133477353136293"
$ws.Range("C10").Value = [double]"0.02297144592282518"
$ws.Range("B11").Value = [double]"-0.19568749936007"
$ws.Range("C11").Value = [double]"0.006268444883292808"
$ws.Range("B12").Value = [double]"0.7712764181279149"
$ws.Range("C12").Value = [double]"6.918746823024867e-13"
$ws.Range("B13").Value = [double]"0.3068252796891577"
$ws.Range("C13").Value = [double]"0.0267019718485558"
$ws.Range("B14").Value = [double]"-3.621887733156215e-05"
$ws.Range("C14").Value = [double]"0.02038080667909046"
$ws.Range("B15").Value = [double]"-2.506817598779832e-08"
$ws.Range("C15").Value = [double]"0.1142897276652645"
$ws.Range("B16").Value = [double]"0.001168203016942514"
$ws.Range("C16").Value = [double]"0.9631831162568707"
$ws.Range("B17").Value = [double]"0.08399460100099959"
$ws.Range("C17").Value = [double]"6.859583663380961e-07"
$ws.Range("B18").Value = [double]"-1.375296981638557"
$ws.Range("C18").Value = [double]"1.530577196407345e-05"
$ws.Range("B19").Value = [double]"-0.0106764399941918"
$ws.Range("C19").Value = [double]"0.007510706155539595"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.002775306013454746"
$ws.Range("C20").Value = [double]"0.243688903418413"
$ws.Range("B21").Value = [double]"0.3407404490863728"
$ws.Range("C21").Value = [double]"0.3544971224137421"
$ws.Range("B22").Value = [double]"0.05623511295582849"
$ws.Range("C22").Value = [double]"0.8627339284275882"
$ws.Rows.Item(23).Delete()

# --- Update data for sheet (physical position 9): summ5 -> summ0 ---
$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = [double]"-0.9278851217479852"
$ws.Range("C2").Value = [double]"0.5130861765126264"
$ws.Range("B3").Value = [double]"1.448164257402054"
$ws.Range("C3").Value = [double]"0.283830985550421"
$ws.Range("B4").Value = [double]"1.562299395734232"
$ws.Range("C4").Value = [double]"0.2503524233673374"
$ws.Range("B5").Value = [double]"0.1886359828496159"
$ws.Range("C5").Value = [double]"0.8891515736204705"
$ws.Range("B6").Value = [double]"0.4965964798122867"
$ws.Range("C6").Value = [double]"0.7137644645988959"
$ws.Range("B7").Value = [double]"0.8001146428708004"
$ws.Range("C7").Value = [double]"0.5560207759731045"
$ws.Range("B8").Value = [double]"-0.03576425398864534"
$ws.Range("C8").Value = [double]"0.6179464340662053"
$ws.Range("B9").Value = [double]"0.0004868006131975214"
$ws.Range("C9").Value = [double]"8.545356079997052e-59"
$ws.Range("B10").Value = [double]"0.006430450655739865"
$ws.Range("C10").Value = [double]"0.01742443454492013"
$ws.Range("B11").Value = [double]"-0.112377653143779"
$ws.Range("C11").Value = [double]"0.1145164146710863"
$ws.Range("B12").Value = [double]"0.7749912904964009"
$ws.Range("C12").Value = [double]"5.268240631013593e-13"
$ws.Range("B13").Value = [double]"0.3386574319353107"
$ws.Range("C13").Value = [double]"0.01488396105932946"
$ws.Range("B14").Value = [double]"-4.356678811377666e-05"
$ws.Range("C14").Value = [double]"0.005118195845825989"
$ws.Range("B15").Value = [double]"-1.497312725482666e-08"
$ws.Range("C15").Value = [double]"0.3450190810519501"
$ws.Range("B16").Value = [double]"0.004288376752578158"
$ws.Range("C16").Value = [double]"0.864657277860603"
$ws.Range("B17").Value = [double]"0.08783744284477708"
$ws.Range("C17").Value = [double]"2.314154223603819e-07"
$ws.Range("B18").Value = [double]"-1.682543718511375"
$ws.Range("C18").Value = [double]"1.066319124183302e-07"
$ws.Range("B19").Value = [double]"-0.008344702594026124"
$ws.Range("C19").Value = [double]"0.03619207346793935"
$ws.Range("A20").Value = "street_length"
$ws.Range("B20").Value = [double]"-0.00361065350548612"
$ws.Range("C20").Value = [double]"0.1263681407258456"
$ws.Range("B21").Value = [double]"0.2061999233882346"
$ws.Range("C21").Value = [double]"0.5754728244946443"
$ws.Range("B22").Value = [double]"-0.06945734914697155"
$ws.Range("C22").Value = [double]"0.8326319144612984"
$ws.Rows.Item(23).Delete()

